# Key based equity data for segmentation.
# Adds three new SEPTA bus routes (59, 66, 75) into the sorted equity table,
# cloning the nearest route's metrics as placeholder values (matching the
# source workbook's edit), and updates the view state to where the author
# was last working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert route 59 (between route 58 and route 60) -----------------------
$ws.Rows("52:52").Insert()
$ws.Cells.Item(52, 1).Value = 59
$ws.Cells.Item(52, 2).Value = 29229.474381905398
$ws.Cells.Item(52, 3).Value = 0.25942586657799699
$ws.Cells.Item(52, 4).Value = 0.46790691798670903
$ws.Cells.Item(52, 5).Value = 0.26481719999901099
$ws.Cells.Item(52, 6).Value = 0.13589523821621899
$ws.Cells.Item(52, 7).Value = 0.239267174407649
$ws.Cells.Item(52, 8).Value = 0.19739211695040701
$ws.Cells.Item(52, 9).Value = 0.019141931902294601

# --- Insert route 66 (between route 65 and route 67) ------------------------
$ws.Rows("58:58").Insert()
$ws.Cells.Item(58, 1).Value = 66
$ws.Cells.Item(58, 2).Value = 29229.474381905398
$ws.Cells.Item(58, 3).Value = 0.25942586657799699
$ws.Cells.Item(58, 4).Value = 0.46790691798670903
$ws.Cells.Item(58, 5).Value = 0.26481719999901099
$ws.Cells.Item(58, 6).Value = 0.13589523821621899
$ws.Cells.Item(58, 7).Value = 0.239267174407649
$ws.Cells.Item(58, 8).Value = 0.19739211695040701
$ws.Cells.Item(58, 9).Value = 0.019141931902294601

# --- Insert route 75 (between route 73 and route 77) ------------------------
$ws.Rows("63:63").Insert()
$ws.Cells.Item(63, 1).Value = 75
$ws.Cells.Item(63, 2).Value = 24303.769836578402
$ws.Cells.Item(63, 3).Value = 0.45287635394213799
$ws.Cells.Item(63, 4).Value = 0.14213104353922101
$ws.Cells.Item(63, 5).Value = 0.67426043590933904
$ws.Cells.Item(63, 6).Value = 0.17263067793999901
$ws.Cells.Item(63, 7).Value = 0.489576158569835
$ws.Cells.Item(63, 8).Value = 0.186064229550646
$ws.Cells.Item(63, 9).Value = 0.039203709356546601

# --- Restore the author's scroll position / active cell ---------------------
$ws.Range("C46").Select()
$excel.ActiveWindow.ScrollRow = 30
